# xlsx boolean parsing fix
#
# Adds a new column E to Sheet1 exercising TRUE()/FALSE() formula cells
# (boolean parsing test data) together with a new shared string "e" for
# its header, and restores C3 (a pre-existing TRUE() formula cell whose
# style had accidentally been pointed at a date/time number format) back
# to the General format now that boolean parsing is correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell for column E
$ws.Range("E1").Value = "e"

# New boolean-formula cells in column E (exercise TRUE()/FALSE() parsing)
$ws.Range("E2").Formula = "=FALSE()"
$ws.Range("E3").Formula = "=TRUE()"
$ws.Range("E4").Formula = "=TRUE()"

# Give the new boolean cells plain numeric formats (mirrors the other
# numeric columns on the sheet)
$ws.Range("E2:E3").NumberFormat = "0"
$ws.Range("E4").NumberFormat = "D/M/YYYY H:MM"

# C3 previously carried a leftover date/time style; now that TRUE()/FALSE()
# formulas parse correctly it is reset back to the General format
$ws.Range("C3").NumberFormat = "General"

# Widen the new column to fit its content
$ws.Columns.Item(5).ColumnWidth = 10.87

# Move the active selection to the newly added cell
[void]$ws.Range("E4").Select()

# Cosmetic: shrink the sheet-tab-bar split ratio
$excel.ActiveWindow.TabRatio = 500
